# [web] clickIfPresent(locator): NEW command that clicks on an element only
# if it can be found.
#
# The '#system' sheet holds the backing lists for the data-validation
# dropdowns used on the 'MacroLibrary' sheet. Column AE (rows 2-158) backs
# the named range "web" - the alphabetically sorted list of WEB command
# signatures. We need to insert the new entry "clickIfPresent(locator)" in
# its correct alphabetical slot (row 63, directly before "clickOffset(...)"),
# shifting the existing entries at rows 63-158 down to rows 64-159, and
# extend the "web" named range so it covers the new row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("#system")

$newEntry = "clickIfPresent(locator)"
$insertRow = 63
$lastRow = 158
$col = 31    # column AE

# Shift existing values in AE63:AE158 down by one row (AE64:AE159).
# Walk bottom-up so we never clobber a value before it has been copied.
for ($r = $lastRow; $r -ge $insertRow; $r--) {
    $v = $ws.Cells.Item($r, $col).Value2
    $ws.Cells.Item($r + 1, $col).Value = $v
}

# Place the new command where it belongs alphabetically.
$ws.Cells.Item($insertRow, $col).Value = $newEntry

# Grow the "web" named range so the dropdown picks up the new last row.
$webName = $wb.Names.Item("web")
$webName.RefersTo = '=''#system''!$AE$2:$AE$159'
